$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 -> becomes old row 4's data
$ws.Range("A3").Value = 130826784
$ws.Range("B3").Value = 57884
$ws.Range("E3").Value = 100109
$ws.Range("F3").Value = "Tretåig hackspett"
$ws.Range("G3").Value = "Picoides tridactylus"
$ws.Range("M3").Value = "färska spår"
$ws.Range("P3").Value = "Brännan, Kälom, Offerdal, Jmt"
$ws.Range("Q3").Value = 461233
$ws.Range("R3").Value = 7039438
$ws.Range("S3").Value = 10
$ws.Range("Z3").Value = "11:37"
$ws.Range("AB3").Value = "11:37"
$ws.Range("AC3").Value = "Födosök barkfläk"

# Row 4 -> becomes old row 5's data
$ws.Range("A4").Value = 130825852
$ws.Range("B4").Value = 57884
$ws.Range("E4").Value = 100109
$ws.Range("F4").Value = "Tretåig hackspett"
$ws.Range("G4").Value = "Picoides tridactylus"
$ws.Range("M4").Value = "färska spår"
$ws.Range("P4").Value = "Flinktorpet, Kälom, Offerdal, Jmt"
$ws.Range("Q4").Value = 460952
$ws.Range("R4").Value = 7039723
$ws.Range("S4").Value = 15
$ws.Range("Z4").Value = "10:42"
$ws.Range("AB4").Value = "10:42"
$ws.Range("AC4").Value = "Barkfläkta grövre och klenare granar."

# Row 5 -> becomes old row 3's data
$ws.Range("A5").Value = 130825823
$ws.Range("B5").Value = 57881
$ws.Range("E5").Value = 100049
$ws.Range("F5").Value = "Spillkråka"
$ws.Range("G5").Value = "Dryocopus martius"
$ws.Range("M5").Value = "äldre spår"
$ws.Range("P5").Value = "Flinktorpet, Kälom, Offerdal, Jmt"
$ws.Range("Q5").Value = 460947
$ws.Range("R5").Value = 7039711
$ws.Range("S5").Value = 10
$ws.Range("Z5").Value = "10:38"
$ws.Range("AB5").Value = "10:38"
$ws.Range("AC5").Value = "Födosökshål på äldre döende gran."
